$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled/recalculated data
$ws.Range("F2").Value = 1
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = -1
